# "cambios en la grafica"
# Duplicate a block of previously-seen rows (with a couple of new ones mixed
# in) onto the bottom of the log, same as rows 45-52 were appended earlier:
# column A keeps the hyperlink look (blue/underline) and a real hyperlink
# pointing at the local file, while the "sae" column keeps its numeric-looking
# codes stored as text (leading apostrophe), matching the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ n=53; a="/Users/alexisjankowicz/Downloads/01 DETENIDO-ART 104-  P PATRICIOS 41 .docx.pdf"; b="28  de  agosto  de  2025"; c="16:08"; d="Art  104"; e="45965420"; f="Inspector  LP  4532  Aguirre  Alan"; g="Oficial   lp  29049  Guzamn  Yesica"; h="Personal  Contratado  Lamboglia  Lautaro"; i="UN  (  1)DETENIDO" },
  @{ n=54; a="/Users/alexisjankowicz/Downloads/CABALLITOO109 - UN (01) DETENIDO.pdf"; b="30 de agosto de 2025"; c="03:29"; d="ROBO"; e="S/D"; f="Inspector ARGUELLO Leandro"; g="Oficial Primero ALCARAZ Ariel"; h="Aux ROJAS Priscila"; i="UN (01) DETENIDO" },
  @{ n=55; a="/Users/alexisjankowicz/Downloads/CONSTITUCIÓN11 - DOS (02) DETENIDOS.pdf"; b="30 de agosto de 2025"; c="03:48"; d="ROBO"; e="45981848"; f="Inspector ARGUELLO Leandro"; g="Oficial Primero ALCARAZ Ariel"; h="Aux ROJAS Priscila"; i="DOS (02) DETENIDOS" },
  @{ n=56; a="/Users/alexisjankowicz/Downloads/RECOLETAO115 - UN (01) DETENIDO.pdf"; b="29 de agosto de 2025"; c="23:28"; d="AV INCENDIO"; e="S/D"; f="Inspector ARGUELLO Leandro"; g="Oficial Primero VILLALBA Nicolas"; h="Aux GARRIGO Agustina"; i="UN (01) DETENIDO" },
  @{ n=57; a="/Users/alexisjankowicz/Downloads/SAN NICOLAS 67 - UN (01) DETENIDO.pdf"; b="29 de agosto de 2025"; c="22:01"; d="TTVA HURTO"; e="45979944"; f="Inspector ARGUELLO Leandro"; g="Inspector RODRIGUEZ Andres"; h="Cont DONADON Nicolas"; i="UN (01) DETENIDO" },
  @{ n=58; a="/Users/alexisjankowicz/Downloads/01 DETENIDO-ART 104-  P PATRICIOS 41 .docx.pdf"; b="28  de  agosto  de  2025"; c="16:08"; d="Art  104"; e="45965420"; f="Inspector  LP  4532  Aguirre  Alan"; g="Oficial   lp  29049  Guzamn  Yesica"; h="Personal  Contratado  Lamboglia  Lautaro"; i="UN  (  1)DETENIDO" },
  @{ n=59; a="/Users/alexisjankowicz/Downloads/CABALLITOO109 - UN (01) DETENIDO.pdf"; b="30 de agosto de 2025"; c="03:29"; d="ROBO"; e="S/D"; f="Inspector ARGUELLO Leandro"; g="Oficial Primero ALCARAZ Ariel"; h="Aux ROJAS Priscila"; i="UN (01) DETENIDO" },
  @{ n=60; a="/Users/alexisjankowicz/Downloads/CONSTITUCIÓN11 - DOS (02) DETENIDOS.pdf"; b="30 de agosto de 2025"; c="03:48"; d="ROBO"; e="45981848"; f="Inspector ARGUELLO Leandro"; g="Oficial Primero ALCARAZ Ariel"; h="Aux ROJAS Priscila"; i="DOS (02) DETENIDOS" },
  @{ n=61; a="/Users/alexisjankowicz/Downloads/RECOLETAO115 - UN (01) DETENIDO.pdf"; b="29 de agosto de 2025"; c="23:28"; d="AV INCENDIO"; e="S/D"; f="Inspector ARGUELLO Leandro"; g="Oficial Primero VILLALBA Nicolas"; h="Aux GARRIGO Agustina"; i="UN (01) DETENIDO" },
  @{ n=62; a="/Users/alexisjankowicz/Downloads/SAN NICOLAS 67 - UN (01) DETENIDO.pdf"; b="29 de agosto de 2025"; c="22:01"; d="TTVA HURTO"; e="45979944"; f="Inspector ARGUELLO Leandro"; g="Inspector RODRIGUEZ Andres"; h="Cont DONADON Nicolas"; i="UN (01) DETENIDO" },
  @{ n=63; a="/Users/alexisjankowicz/Downloads/1 CONTRAVENTOR- SAE 45976643.pdf"; b="29 de Agosto del 2025"; c="17:00"; d="ARTICULO 91 Y 239"; e="45976643"; f="Oficial Mayor MARTINEZ, Ricardo"; g="Oficial PEREZ, Ramón"; h="SEGOVIA, Noelia"; i="UN (01) CONTRAVENTOR" },
  @{ n=64; a="/Users/alexisjankowicz/Downloads/CABALLITOO109 - UN (01) DETENIDO.pdf"; b="30 de agosto de 2025"; c="03:29"; d="ROBO"; e="S/D"; f="Inspector ARGUELLO Leandro"; g="Oficial Primero ALCARAZ Ariel"; h="Aux ROJAS Priscila"; i="UN (01) DETENIDO" }
)

foreach ($row in $rows) {
  $n = $row.n

  # Column A: file path text + hyperlink styling (blue, underlined) + the
  # actual external hyperlink relationship, same as the existing A45:A52 block.
  # Add the hyperlink first, then pin down the font explicitly -- Add() alone
  # pulls in the theme's "Hyperlink" look, and we want the sheet's own
  # plain blue/underline font instead (same as columns A2:A52).
  $cellA = $ws.Cells.Item($n, 1)
  $cellA.Value = $row.a
  $ws.Hyperlinks.Add($cellA, $row.a)
  $cellA.Font.Underline = 2
  $cellA.Font.Color = 16711680

  $ws.Cells.Item($n, 2).Value = $row.b
  $ws.Cells.Item($n, 3).Value = $row.c
  $ws.Cells.Item($n, 4).Value = $row.d

  # Column E ("sae"): some values are numeric-looking codes that must stay as
  # text (leading apostrophe forces text, like the earlier rows in the file);
  # others ("S/D") are already plain text.
  $cellE = $ws.Cells.Item($n, 5)
  if ($row.e -match '^\d+$') {
    $cellE.Value = "'" + $row.e
  } else {
    $cellE.Value = $row.e
  }

  $ws.Cells.Item($n, 6).Value = $row.f
  $ws.Cells.Item($n, 7).Value = $row.g
  $ws.Cells.Item($n, 8).Value = $row.h
  $ws.Cells.Item($n, 9).Value = $row.i
}
